$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coin price/volume refresh (cryptos list update).
# Column D ("Price") holds plain-text numbers such as "26.388.65" or "1.001";
# format the cells as Text first so Excel does not reinterpret/round the literal
# string (e.g. dropping a trailing zero or treating dotted values as numbers).

$priceUpdates = @(
    @{ Cell = "D2"; Value = '26.388.65' },
    @{ Cell = "D3"; Value = '1.843.09' },
    @{ Cell = "D4"; Value = '1.001' },
    @{ Cell = "D5"; Value = '259.97' },
    @{ Cell = "D7"; Value = '0.5116' },
    @{ Cell = "D8"; Value = '0.3215' },
    @{ Cell = "D9"; Value = '0.06745' },
    @{ Cell = "D10"; Value = '19.15' },
    @{ Cell = "D11"; Value = '0.7718' },
    @{ Cell = "D12"; Value = '0.07690' },
    @{ Cell = "D13"; Value = '1.841.43' },
    @{ Cell = "D15"; Value = '5.028' },
    @{ Cell = "D16"; Value = '1.000' },
    @{ Cell = "D17"; Value = '14.09' },
    @{ Cell = "D19"; Value = '0.000007894' },
    @{ Cell = "D20"; Value = '26.454.77' },
    @{ Cell = "D21"; Value = '2.096.62' },
    @{ Cell = "D22"; Value = '4.586' },
    @{ Cell = "D23"; Value = '9.546' },
    @{ Cell = "D24"; Value = '5.984' },
    @{ Cell = "D25"; Value = '2.349' },
    @{ Cell = "D26"; Value = '144.99' },
    @{ Cell = "D27"; Value = '1.653' },
    @{ Cell = "D28"; Value = '16.94' },
    @{ Cell = "D29"; Value = '110.90' },
    @{ Cell = "D30"; Value = '4.199' },
    @{ Cell = "D31"; Value = '4.170' },
    @{ Cell = "D32"; Value = '0.08707' },
    @{ Cell = "D33"; Value = '0.04833' },
    @{ Cell = "D35"; Value = '2.838' },
    @{ Cell = "D36"; Value = '0.6899' },
    @{ Cell = "D37"; Value = '3.086' },
    @{ Cell = "D38"; Value = '0.01811' },
    @{ Cell = "D39"; Value = '2.218' },
    @{ Cell = "D40"; Value = '0.4920' },
    @{ Cell = "D41"; Value = '113.45' },
    @{ Cell = "D42"; Value = '0.9053' },
    @{ Cell = "D43"; Value = '6.112' },
    @{ Cell = "D45"; Value = '7.803' },
    @{ Cell = "D46"; Value = '0.4261' },
    @{ Cell = "D47"; Value = '0.1283' },
    @{ Cell = "D48"; Value = '9.139' },
    @{ Cell = "D49"; Value = '0.05900' },
    @{ Cell = "D50"; Value = '35.20' },
    @{ Cell = "D51"; Value = '1.435' }
)

$volumeUpdates = @(
    @{ Cell = "E2"; Value = '  +0.54%  ' },
    @{ Cell = "E3"; Value = '  -0.55%  ' },
    @{ Cell = "E4"; Value = '  +0.36%  ' },
    @{ Cell = "E5"; Value = '  -7.18%  ' },
    @{ Cell = "E6"; Value = '  +0.24%  ' },
    @{ Cell = "E7"; Value = '  -0.33%  ' },
    @{ Cell = "E8"; Value = '  -8.49%  ' },
    @{ Cell = "E9"; Value = '  -1.57%  ' },
    @{ Cell = "E10"; Value = '  -4.41%  ' },
    @{ Cell = "E11"; Value = '  -4.78%  ' },
    @{ Cell = "E12"; Value = '  -1.00%  ' },
    @{ Cell = "E13"; Value = '  -0.45%  ' },
    @{ Cell = "E14"; Value = '  -0.98%  ' },
    @{ Cell = "E15"; Value = '  -1.54%  ' },
    @{ Cell = "E16"; Value = '  +0.28%  ' },
    @{ Cell = "E17"; Value = '  -0.97%  ' },
    @{ Cell = "E19"; Value = '  -2.57%  ' },
    @{ Cell = "E20"; Value = '  +0.78%  ' },
    @{ Cell = "E21"; Value = '  +0.82%  ' },
    @{ Cell = "E22"; Value = '  -4.24%  ' },
    @{ Cell = "E23"; Value = '  -5.46%  ' },
    @{ Cell = "E24"; Value = '  -3.85%  ' },
    @{ Cell = "E25"; Value = '  -1.26%  ' },
    @{ Cell = "E26"; Value = '  +0.29%  ' },
    @{ Cell = "E27"; Value = '  -0.73%  ' },
    @{ Cell = "E28"; Value = '  -1.92%  ' },
    @{ Cell = "E29"; Value = '  +0.59%  ' },
    @{ Cell = "E30"; Value = '  -4.05%  ' },
    @{ Cell = "E31"; Value = '  -3.59%  ' },
    @{ Cell = "E32"; Value = '  -0.81%  ' },
    @{ Cell = "E33"; Value = '  -1.79%  ' },
    @{ Cell = "E34"; Value = '  -3.85%  ' },
    @{ Cell = "E35"; Value = '  -0.09%  ' },
    @{ Cell = "E36"; Value = '  -7.06%  ' },
    @{ Cell = "E37"; Value = '  -4.60%  ' },
    @{ Cell = "E38"; Value = '  -2.54%  ' },
    @{ Cell = "E39"; Value = '  -7.75%  ' },
    @{ Cell = "E40"; Value = '  -5.03%  ' },
    @{ Cell = "E41"; Value = '  -2.43%  ' },
    @{ Cell = "E42"; Value = '  -6.27%  ' },
    @{ Cell = "E43"; Value = '  -2.45%  ' },
    @{ Cell = "E44"; Value = '  +0.18%  ' },
    @{ Cell = "E45"; Value = '  -2.96%  ' },
    @{ Cell = "E46"; Value = '  -6.11%  ' },
    @{ Cell = "E47"; Value = '  -5.72%  ' },
    @{ Cell = "E48"; Value = '  -2.97%  ' },
    @{ Cell = "E49"; Value = '  -0.41%  ' },
    @{ Cell = "E50"; Value = '  -3.23%  ' },
    @{ Cell = "E51"; Value = '  -4.66%  ' }
)

foreach ($u in $priceUpdates) {
    $ws.Range($u.Cell).NumberFormat = "@"
}
foreach ($u in $priceUpdates) {
    $ws.Range($u.Cell).Value = $u.Value
}
foreach ($u in $volumeUpdates) {
    $ws.Range($u.Cell).Value = $u.Value
}
